$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Q1_20_21": the F9 row (row 4) is replaced by the Columbia
# row (old row 5), and the old row 5 is removed -> B2:K4
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Q1_20_21")

$ws1.Range("B4").Value = "Columbia"
$ws1.Range("C4").ClearContents()
$ws1.Range("D4").Value = 0.38
$ws1.Range("E4").Value = 0.63
$ws1.Range("F4").Value = "Poor"
$ws1.Range("G4").ClearContents()
$ws1.Range("H4").ClearContents()
$ws1.Range("I4").Value = 1172
$ws1.Range("J4").Value = 738.36

$ws1.Rows("5:5").Delete()

# ---------------------------------------------------------------
# Sheet "Q4_19_20": rows 4 (A13) and 5 (Columbia) are removed, so
# the old row 6 (F9) becomes the new row 4 -> B2:K4
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Q4_19_20")

$ws2.Rows("4:5").Delete()

# ---------------------------------------------------------------
# Sheet "Count": update totals/counts that changed as a result of
# the data restructuring above
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Count")

$ws3.Range("D4").Value = 0
$ws3.Range("C6").Value = 0
$ws3.Range("D7").Value = 928
$ws3.Range("C11").Value = 1172
$ws3.Range("D11").Value = 3759

$ws3.Range("D16").Value = 0
$ws3.Range("C18").Value = 0
$ws3.Range("D19").Value = 1
$ws3.Range("C23").Value = 2
$ws3.Range("D23").Value = 2
